$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.884.15'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '2.462.29'
$ws.Range("E3").Value = '  -0.63%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'571.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.12%  '

$ws.Range("D6").Value = "'146.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -1.83%  '

$ws.Range("D9").Value = "'0.111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.17%  '

$ws.Range("E11").Value = '  -1.81%  '

$ws.Range("E12").Value = '  -1.80%  '

$ws.Range("D13").Value = "'28.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("E14").Value = '  -3.07%  '

$ws.Range("D15").Value = '2.900.75'
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("D16").Value = '62.747.10'
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("D17").Value = '2.460.40'
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("D18").Value = "'7.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.31%  '

$ws.Range("D19").Value = "'10.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.04%  '

$ws.Range("E20").Value = '  -1.33%  '

$ws.Range("D21").Value = "'322.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.19%  '

$ws.Range("D22").Value = "'4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").Value = "'10.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.90%  '

$ws.Range("D25").Value = "'64.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.30%  '

$ws.Range("D26").Value = "'647.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.04%  '

$ws.Range("D27").Value = '2.580.12'
$ws.Range("E27").Value = '  -0.95%  '

$ws.Range("D28").Value = '0.0₃0961'
$ws.Range("E28").Value = '  -3.78%  '

$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("E30").Value = '  -3.21%  '

$ws.Range("D31").Value = "'7.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.55%  '

$ws.Range("E32").Value = '  -3.04%  '

$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = "'1.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.55%  '

$ws.Range("D36").Value = "'4.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.98%  '

$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("D38").Value = "'18.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.50%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = "'5.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.57%  '

$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = "'0.364"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.28%  '

$ws.Range("D41").Value = "'2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.51%  '

$ws.Range("D42").Value = "'1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.20%  '

$ws.Range("D43").Value = '0.0₆0312'
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D45").Value = "'152.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.85%  '

$ws.Range("D46").Value = "'15.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.86%  '

$ws.Range("D47").Value = "'3.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.12%  '

$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("D49").Value = "'20.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.49%  '

$ws.Range("D50").Value = "'0.0506"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("D51").Value = "'0.0902"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.84%  '
